$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 502 (weekly update), shifting existing
# rows 502:565 down to 503:566.
$ws.Rows.Item(502).Insert()

# Populate the newly inserted row with this week's data.
$ws.Cells.Item(502, 1).Value = 8
$ws.Cells.Item(502, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(502, 3).Value = "Coquimbo"
$ws.Cells.Item(502, 4).Value = 45124
$ws.Cells.Item(502, 5).Value = 4
$ws.Cells.Item(502, 6).Value = 100114013
$ws.Cells.Item(502, 7).Value = "Zanahoria"
$ws.Cells.Item(502, 8).Value = "Sin especificar"
$ws.Cells.Item(502, 9).Value = "Primera"
$ws.Cells.Item(502, 10).Value = 540
$ws.Cells.Item(502, 11).Value = 5800
$ws.Cells.Item(502, 12).Value = 6000
$ws.Cells.Item(502, 13).Value = 5900
$ws.Cells.Item(502, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(502, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(502, 16).Value = 295
$ws.Cells.Item(502, 17).Value = 20
$ws.Cells.Item(502, 18).Value = "Hortaliza"
